$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "TODO's" sheet: status of "add random procedure generated dungeons" (row 6)
# moves from "todo" to "in-progress"
# ---------------------------------------------------------------------------
$wsTodos = $wb.Worksheets.Item("TODO's")
$wsTodos.Range("C6").Value = "in-progress"

# update the selection on that sheet (was B10) without leaving it as the
# active tab - we restore the "Logs" sheet as active afterwards.
$wsTodos.Range("C7").Select()

# ---------------------------------------------------------------------------
# "Logs" sheet: add four new dev-log entries (rows 47-50)
# ---------------------------------------------------------------------------
$wsLogs = $wb.Worksheets.Item("Logs")

# Row 47
$wsLogs.Range("A46:B46").Copy()
$wsLogs.Range("A47:B47").PasteSpecial(-4122)
$wsLogs.Range("A47").Value = 45483
$wsLogs.Range("B47").Value = "fix existing generation of exit. Improve to reduce complexity, first attempts to add additional room"

# Row 48
$wsLogs.Range("A47:B47").Copy()
$wsLogs.Range("A48:B48").PasteSpecial(-4122)
$wsLogs.Range("A48").Value = 45485
$wsLogs.Range("B48").Value = "DONE random floor and exist generated perfectly! Next step to clean up that generator class and add the walls "

# Row 49
$wsLogs.Range("A48:B48").Copy()
$wsLogs.Range("A49:B49").PasteSpecial(-4122)
$wsLogs.Range("A49").Value = 45492
$wsLogs.Range("B49").Value = "Added system of room saving - and now each generated room save some info. Next step to use it while walls generation - to catch where wall should be and where not."

# Row 50
$wsLogs.Range("A49:B49").Copy()
$wsLogs.Range("A50:B50").PasteSpecial(-4122)
$wsLogs.Range("A50").Value = 45493
$wsLogs.Range("B50").Value = "Add deadend rooms and modify the common room schema"

# move selection to the new last row and keep "Logs" as the active sheet
$wsLogs.Range("B51").Select()
